$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-21 Friday", 2) | Out-Null
$d.Content.Find.Execute("68×87=5916", $true, $false, $false, $false, $false, $true, 1, $false, "63×30=1890", 2) | Out-Null
$d.Content.Find.Execute("38×21=798", $true, $false, $false, $false, $false, $true, 1, $false, "34×57=1938", 2) | Out-Null
$d.Content.Find.Execute("42×12=504", $true, $false, $false, $false, $false, $true, 1, $false, "20×97=1940", 2) | Out-Null
$d.Content.Find.Execute("21×48=1008", $true, $false, $false, $false, $false, $true, 1, $false, "89×70=6230", 2) | Out-Null
$d.Content.Find.Execute("55×26=1430", $true, $false, $false, $false, $false, $true, 1, $false, "56×20=1120", 2) | Out-Null
$d.Content.Find.Execute("89×75=6675", $true, $false, $false, $false, $false, $true, 1, $false, "61×79=4819", 2) | Out-Null
$d.Content.Find.Execute("62×57=3534", $true, $false, $false, $false, $false, $true, 1, $false, "58×93=5394", 2) | Out-Null
$d.Content.Find.Execute("93×90=8370", $true, $false, $false, $false, $false, $true, 1, $false, "24×97=2328", 2) | Out-Null
$d.Content.Find.Execute("98×23=2254", $true, $false, $false, $false, $false, $true, 1, $false, "77×24=1848", 2) | Out-Null
$d.Content.Find.Execute("38×71=2698", $true, $false, $false, $false, $false, $true, 1, $false, "18×13=234", 2) | Out-Null
$d.Content.Find.Execute("76×20=1520", $true, $false, $false, $false, $false, $true, 1, $false, "64×73=4672", 2) | Out-Null
$d.Content.Find.Execute("76×90=6840", $true, $false, $false, $false, $false, $true, 1, $false, "18×91=1638", 2) | Out-Null
$d.Content.Find.Execute("48×57=2736", $true, $false, $false, $false, $false, $true, 1, $false, "87×19=1653", 2) | Out-Null
$d.Content.Find.Execute("53×57=3021", $true, $false, $false, $false, $false, $true, 1, $false, "89×90=8010", 2) | Out-Null
$d.Content.Find.Execute("14×75=1050", $true, $false, $false, $false, $false, $true, 1, $false, "62×89=5518", 2) | Out-Null
$d.Content.Find.Execute("71×65=4615", $true, $false, $false, $false, $false, $true, 1, $false, "97×45=4365", 2) | Out-Null
$d.Content.Find.Execute("15×55=825", $true, $false, $false, $false, $false, $true, 1, $false, "19×26=494", 2) | Out-Null
$d.Content.Find.Execute("52×72=3744", $true, $false, $false, $false, $false, $true, 1, $false, "100×63=6300", 2) | Out-Null
$d.Content.Find.Execute("62×19=1178", $true, $false, $false, $false, $false, $true, 1, $false, "63×16=1008", 2) | Out-Null
$d.Content.Find.Execute("57×62=3534", $true, $false, $false, $false, $false, $true, 1, $false, "84×74=6216", 2) | Out-Null
$d.Content.Find.Execute("19×53=1007", $true, $false, $false, $false, $false, $true, 1, $false, "67×65=4355", 2) | Out-Null
$d.Content.Find.Execute("34×83=2822", $true, $false, $false, $false, $false, $true, 1, $false, "44×20=880", 2) | Out-Null
$d.Content.Find.Execute("37×13=481", $true, $false, $false, $false, $false, $true, 1, $false, "57×10=570", 2) | Out-Null
$d.Content.Find.Execute("25×42=1050", $true, $false, $false, $false, $false, $true, 1, $false, "53×56=2968", 2) | Out-Null
$d.Content.Find.Execute("67×49=3283", $true, $false, $false, $false, $false, $true, 1, $false, "94×55=5170", 2) | Out-Null
$d.Content.Find.Execute("66×76=5016", $true, $false, $false, $false, $false, $true, 1, $false, "45×30=1350", 2) | Out-Null
$d.Content.Find.Execute("55×67=3685", $true, $false, $false, $false, $false, $true, 1, $false, "58×72=4176", 2) | Out-Null
$d.Content.Find.Execute("67×99=6633", $true, $false, $false, $false, $false, $true, 1, $false, "22×83=1826", 2) | Out-Null
$d.Content.Find.Execute("17×81=1377", $true, $false, $false, $false, $false, $true, 1, $false, "99×67=6633", 2) | Out-Null
$d.Content.Find.Execute("22×23=506", $true, $false, $false, $false, $false, $true, 1, $false, "82×26=2132", 2) | Out-Null
$d.Content.Find.Execute("42×65=2730", $true, $false, $false, $false, $false, $true, 1, $false, "89×54=4806", 2) | Out-Null
$d.Content.Find.Execute("18×99=1782", $true, $false, $false, $false, $false, $true, 1, $false, "35×15=525", 2) | Out-Null
$d.Content.Find.Execute("94×57=5358", $true, $false, $false, $false, $false, $true, 1, $false, "51×24=1224", 2) | Out-Null
$d.Content.Find.Execute("39×92=3588", $true, $false, $false, $false, $false, $true, 1, $false, "59×35=2065", 2) | Out-Null
$d.Content.Find.Execute("73×32=2336", $true, $false, $false, $false, $false, $true, 1, $false, "64×79=5056", 2) | Out-Null
$d.Content.Find.Execute("38×60=2280", $true, $false, $false, $false, $false, $true, 1, $false, "59×84=4956", 2) | Out-Null
$d.Content.Find.Execute("39×10=390", $true, $false, $false, $false, $false, $true, 1, $false, "63×13=819", 2) | Out-Null
$d.Content.Find.Execute("39×89=3471", $true, $false, $false, $false, $false, $true, 1, $false, "50×48=2400", 2) | Out-Null
$d.Content.Find.Execute("80×96=7680", $true, $false, $false, $false, $false, $true, 1, $false, "58×49=2842", 2) | Out-Null
$d.Content.Find.Execute("65×50=3250", $true, $false, $false, $false, $false, $true, 1, $false, "93×13=1209", 2) | Out-Null
$d.Content.Find.Execute("95×25=2375", $true, $false, $false, $false, $false, $true, 1, $false, "59×93=5487", 2) | Out-Null
$d.Content.Find.Execute("43×20=860", $true, $false, $false, $false, $false, $true, 1, $false, "52×60=3120", 2) | Out-Null
$d.Content.Find.Execute("33×100=3300", $true, $false, $false, $false, $false, $true, 1, $false, "93×51=4743", 2) | Out-Null
$d.Content.Find.Execute("29×78=2262", $true, $false, $false, $false, $false, $true, 1, $false, "25×29=725", 2) | Out-Null
$d.Content.Find.Execute("21×46=966", $true, $false, $false, $false, $false, $true, 1, $false, "51×97=4947", 2) | Out-Null
$d.Content.Find.Execute("40×97=3880", $true, $false, $false, $false, $false, $true, 1, $false, "83×13=1079", 2) | Out-Null
$d.Content.Find.Execute("89×35=3115", $true, $false, $false, $false, $false, $true, 1, $false, "55×60=3300", 2) | Out-Null
$d.Content.Find.Execute("76×96=7296", $true, $false, $false, $false, $false, $true, 1, $false, "13×37=481", 2) | Out-Null
$d.Content.Find.Execute("33×22=726", $true, $false, $false, $false, $false, $true, 1, $false, "31×31=961", 2) | Out-Null
$d.Content.Find.Execute("66×29=1914", $true, $false, $false, $false, $false, $true, 1, $false, "54×63=3402", 2) | Out-Null
$d.Content.Find.Execute("46×56=2576", $true, $false, $false, $false, $false, $true, 1, $false, "40×14=560", 2) | Out-Null
$d.Content.Find.Execute("58×62=3596", $true, $false, $false, $false, $false, $true, 1, $false, "65×58=3770", 2) | Out-Null
$d.Content.Find.Execute("88×85=7480", $true, $false, $false, $false, $false, $true, 1, $false, "83×26=2158", 2) | Out-Null
$d.Content.Find.Execute("45×41=1845", $true, $false, $false, $false, $false, $true, 1, $false, "51×91=4641", 2) | Out-Null
$d.Content.Find.Execute("11×94=1034", $true, $false, $false, $false, $false, $true, 1, $false, "35×87=3045", 2) | Out-Null
$d.Content.Find.Execute("37×69=2553", $true, $false, $false, $false, $false, $true, 1, $false, "28×57=1596", 2) | Out-Null
$d.Content.Find.Execute("27×12=324", $true, $false, $false, $false, $false, $true, 1, $false, "69×33=2277", 2) | Out-Null
$d.Content.Find.Execute("57×46=2622", $true, $false, $false, $false, $false, $true, 1, $false, "83×70=5810", 2) | Out-Null
$d.Content.Find.Execute("14×91=1274", $true, $false, $false, $false, $false, $true, 1, $false, "65×88=5720", 2) | Out-Null
$d.Content.Find.Execute("28×60=1680", $true, $false, $false, $false, $false, $true, 1, $false, "45×67=3015", 2) | Out-Null
$d.Content.Find.Execute("100×73=7300", $true, $false, $false, $false, $false, $true, 1, $false, "17×18=306", 2) | Out-Null
$d.Content.Find.Execute("89×98=8722", $true, $false, $false, $false, $false, $true, 1, $false, "97×10=970", 2) | Out-Null
$d.Content.Find.Execute("86×89=7654", $true, $false, $false, $false, $false, $true, 1, $false, "16×19=304", 2) | Out-Null
$d.Content.Find.Execute("32×28=896", $true, $false, $false, $false, $false, $true, 1, $false, "71×47=3337", 2) | Out-Null
$d.Content.Find.Execute("75×93=6975", $true, $false, $false, $false, $false, $true, 1, $false, "70×57=3990", 2) | Out-Null
$d.Content.Find.Execute("78×82=6396", $true, $false, $false, $false, $false, $true, 1, $false, "43×53=2279", 2) | Out-Null
$d.Content.Find.Execute("83×56=4648", $true, $false, $false, $false, $false, $true, 1, $false, "30×39=1170", 2) | Out-Null
$d.Content.Find.Execute("79×99=7821", $true, $false, $false, $false, $false, $true, 1, $false, "59×94=5546", 2) | Out-Null
$d.Content.Find.Execute("51×69=3519", $true, $false, $false, $false, $false, $true, 1, $false, "85×33=2805", 2) | Out-Null
$d.Content.Find.Execute("22×98=2156", $true, $false, $false, $false, $false, $true, 1, $false, "46×65=2990", 2) | Out-Null
$d.Content.Find.Execute("81×79=6399", $true, $false, $false, $false, $false, $true, 1, $false, "30×37=1110", 2) | Out-Null
$d.Content.Find.Execute("31×52=1612", $true, $false, $false, $false, $false, $true, 1, $false, "100×90=9000", 2) | Out-Null
$d.Content.Find.Execute("91×23=2093", $true, $false, $false, $false, $false, $true, 1, $false, "82×26=2132", 2) | Out-Null
$d.Content.Find.Execute("35×95=3325", $true, $false, $false, $false, $false, $true, 1, $false, "87×99=8613", 2) | Out-Null
$d.Content.Find.Execute("42×54=2268", $true, $false, $false, $false, $false, $true, 1, $false, "70×59=4130", 2) | Out-Null
$d.Content.Find.Execute("26×87=2262", $true, $false, $false, $false, $false, $true, 1, $false, "90×12=1080", 2) | Out-Null
$d.Content.Find.Execute("38×11=418", $true, $false, $false, $false, $false, $true, 1, $false, "31×12=372", 2) | Out-Null
$d.Content.Find.Execute("82×78=6396", $true, $false, $false, $false, $false, $true, 1, $false, "77×21=1617", 2) | Out-Null
$d.Content.Find.Execute("18×45=810", $true, $false, $false, $false, $false, $true, 1, $false, "98×15=1470", 2) | Out-Null
$d.Content.Find.Execute("21×16=336", $true, $false, $false, $false, $false, $true, 1, $false, "32×87=2784", 2) | Out-Null
$d.Content.Find.Execute("47×98=4606", $true, $false, $false, $false, $false, $true, 1, $false, "99×28=2772", 2) | Out-Null
$d.Content.Find.Execute("93×92=8556", $true, $false, $false, $false, $false, $true, 1, $false, "93×61=5673", 2) | Out-Null
$d.Content.Find.Execute("47×82=3854", $true, $false, $false, $false, $false, $true, 1, $false, "83×96=7968", 2) | Out-Null
$d.Content.Find.Execute("73×23=1679", $true, $false, $false, $false, $false, $true, 1, $false, "22×64=1408", 2) | Out-Null
$d.Content.Find.Execute("37×42=1554", $true, $false, $false, $false, $false, $true, 1, $false, "20×18=360", 2) | Out-Null
$d.Content.Find.Execute("49×61=2989", $true, $false, $false, $false, $false, $true, 1, $false, "58×68=3944", 2) | Out-Null
$d.Content.Find.Execute("78×99=7722", $true, $false, $false, $false, $false, $true, 1, $false, "62×22=1364", 2) | Out-Null
$d.Content.Find.Execute("67×59=3953", $true, $false, $false, $false, $false, $true, 1, $false, "75×73=5475", 2) | Out-Null
$d.Content.Find.Execute("87×84=7308", $true, $false, $false, $false, $false, $true, 1, $false, "76×74=5624", 2) | Out-Null
$d.Content.Find.Execute("50×87=4350", $true, $false, $false, $false, $false, $true, 1, $false, "60×87=5220", 2) | Out-Null
$d.Content.Find.Execute("24×33=792", $true, $false, $false, $false, $false, $true, 1, $false, "65×71=4615", 2) | Out-Null
$d.Content.Find.Execute("27×58=1566", $true, $false, $false, $false, $false, $true, 1, $false, "99×88=8712", 2) | Out-Null
$d.Content.Find.Execute("53×52=2756", $true, $false, $false, $false, $false, $true, 1, $false, "49×88=4312", 2) | Out-Null
$d.Content.Find.Execute("57×36=2052", $true, $false, $false, $false, $false, $true, 1, $false, "73×69=5037", 2) | Out-Null
$d.Content.Find.Execute("83×61=5063", $true, $false, $false, $false, $false, $true, 1, $false, "41×26=1066", 2) | Out-Null
$d.Content.Find.Execute("66×41=2706", $true, $false, $false, $false, $false, $true, 1, $false, "66×81=5346", 2) | Out-Null
$d.Content.Find.Execute("55×79=4345", $true, $false, $false, $false, $false, $true, 1, $false, "77×25=1925", 2) | Out-Null
$d.Content.Find.Execute("34×68=2312", $true, $false, $false, $false, $false, $true, 1, $false, "83×66=5478", 2) | Out-Null
$d.Content.Find.Execute("23×53=1219", $true, $false, $false, $false, $false, $true, 1, $false, "100×22=2200", 2) | Out-Null
$d.Content.Find.Execute("80×42=3360", $true, $false, $false, $false, $false, $true, 1, $false, "80×90=7200", 2) | Out-Null
